$p = $ppt.ActivePresentation

$oldDate = "9/17/18 3:37 PM"
$newDate = "3/4/19 8:28 PM"

# Update the "Date Placeholder" field text on every slide's Notes Page.
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    $np = $s.NotesPage
    for ($j = 1; $j -le $np.Shapes.Count; $j++) {
        $shp = $np.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Update the Notes Master's Date Placeholder.
$nm = $p.NotesMaster
for ($j = 1; $j -le $nm.Shapes.Count; $j++) {
    $shp = $nm.Shapes.Item($j)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# Update the Handout Master's Date Placeholder.
$hm = $p.HandoutMaster
for ($j = 1; $j -le $hm.Shapes.Count; $j++) {
    $shp = $hm.Shapes.Item($j)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}
